$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows at position 1073, shifting existing data (old rows 1073-1145) down to 1082-1154
$ws.Range("A1073:A1081").EntireRow.Insert()

# Format column B and C of the new rows as Text first so date-like / zero-padded strings are preserved literally
$ws.Range("B1073:C1081").NumberFormat = "@"

$ws.Range("A1073").Value = 1574035200
$ws.Range("B1073").Value = "2019-11-18"
$ws.Range("C1073").Value = "0178"
$ws.Range("D1073").Value = "SEDANIA"
$ws.Range("E1073").Value = 0.165
$ws.Range("F1073").Value = 0.165
$ws.Range("G1073").Value = 0.16
$ws.Range("H1073").Value = 0.16
$ws.Range("I1073").Value = 692900

$ws.Range("A1074").Value = 1574121600
$ws.Range("B1074").Value = "2019-11-19"
$ws.Range("C1074").Value = "0178"
$ws.Range("D1074").Value = "SEDANIA"
$ws.Range("E1074").Value = 0.16
$ws.Range("F1074").Value = 0.16
$ws.Range("G1074").Value = 0.155
$ws.Range("H1074").Value = 0.155
$ws.Range("I1074").Value = 706000

$ws.Range("A1075").Value = 1574208000
$ws.Range("B1075").Value = "2019-11-20"
$ws.Range("C1075").Value = "0178"
$ws.Range("D1075").Value = "SEDANIA"
$ws.Range("E1075").Value = 0.15
$ws.Range("F1075").Value = 0.17
$ws.Range("G1075").Value = 0.15
$ws.Range("H1075").Value = 0.17
$ws.Range("I1075").Value = 1784300

$ws.Range("A1076").Value = 1574294400
$ws.Range("B1076").Value = "2019-11-21"
$ws.Range("C1076").Value = "0178"
$ws.Range("D1076").Value = "SEDANIA"
$ws.Range("E1076").Value = 0.175
$ws.Range("F1076").Value = 0.18
$ws.Range("G1076").Value = 0.17
$ws.Range("H1076").Value = 0.175
$ws.Range("I1076").Value = 4138800

$ws.Range("A1077").Value = 1574380800
$ws.Range("B1077").Value = "2019-11-22"
$ws.Range("C1077").Value = "0178"
$ws.Range("D1077").Value = "SEDANIA"
$ws.Range("E1077").Value = 0.17
$ws.Range("F1077").Value = 0.175
$ws.Range("G1077").Value = 0.165
$ws.Range("H1077").Value = 0.175
$ws.Range("I1077").Value = 468000

$ws.Range("A1078").Value = 1574640000
$ws.Range("B1078").Value = "2019-11-25"
$ws.Range("C1078").Value = "0178"
$ws.Range("D1078").Value = "SEDANIA"
$ws.Range("E1078").Value = 0.175
$ws.Range("F1078").Value = 0.175
$ws.Range("G1078").Value = 0.165
$ws.Range("H1078").Value = 0.17
$ws.Range("I1078").Value = 557300

$ws.Range("A1079").Value = 1574726400
$ws.Range("B1079").Value = "2019-11-26"
$ws.Range("C1079").Value = "0178"
$ws.Range("D1079").Value = "SEDANIA"
$ws.Range("E1079").Value = 0.165
$ws.Range("F1079").Value = 0.175
$ws.Range("G1079").Value = 0.165
$ws.Range("H1079").Value = 0.165
$ws.Range("I1079").Value = 1459100

$ws.Range("A1080").Value = 1574812800
$ws.Range("B1080").Value = "2019-11-27"
$ws.Range("C1080").Value = "0178"
$ws.Range("D1080").Value = "SEDANIA"
$ws.Range("E1080").Value = 0.16
$ws.Range("F1080").Value = 0.16
$ws.Range("G1080").Value = 0.16
$ws.Range("H1080").Value = 0.16
$ws.Range("I1080").Value = 277500

$ws.Range("A1081").Value = 1574899200
$ws.Range("B1081").Value = "2019-11-28"
$ws.Range("C1081").Value = "0178"
$ws.Range("D1081").Value = "SEDANIA"
$ws.Range("E1081").Value = 0.155
$ws.Range("F1081").Value = 0.16
$ws.Range("G1081").Value = 0.155
$ws.Range("H1081").Value = 0.16
$ws.Range("I1081").Value = 264000

# Remove the temporary text formatting from columns B and C so they match the rest of the sheet (no style index)
$ws.Range("B1073:C1081").ClearFormats()
